# Apply the changes described by the diff:
# 1. Rename sheet "High Priority break-up" -> "Interannual update - High Pri"
# 2. Add a new sheet "Major update - High Priority " right after it, containing
#    the original "High Priority break-up" data (the IUCN-only row).
# 3. Update the "Interannual update - High Pri" sheet data: insert a new
#    "Trend New" row before "IUCN", and update the IUCN row's values.
# 4. Update various numeric/text values on other sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet1 "Trends Status": Insufficient Data row totals 235 -> 236 ---
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("B8").Value = 236
$wsTrends.Range("C8").Value = 236

# --- Sheet3 "Priority Status": update species counts ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# --- Sheet4 "Species qualification": rename label, update count ---
$wsSpeciesQual = $wb.Worksheets.Item("Species qualification")
$wsSpeciesQual.Range("A2").Value = "SoIB Assessment"
$wsSpeciesQual.Range("B2").Value = 236

# --- Sheet5 "High Priority break-up": duplicate before editing, so the
#     duplicate preserves the original (unchanged) data and becomes the
#     new "Major update - High Priority " sheet. ---
$wsHighPriority = $wb.Worksheets.Item("High Priority break-up")
$wsHighPriority.Copy($null, $wsHighPriority)
$wsNewCopy = $wb.Worksheets.Item($wsHighPriority.Index + 1)
$wsNewCopy.Name = "Major update - High Priority "

# Rename the original sheet
$wsHighPriority.Name = "Interannual update - High Pri"

# Insert a new row 2 ("Trend New") on the renamed ("Interannual update")
# sheet, pushing the existing "IUCN" row down to row 3, then set values.
$wsHighPriority.Rows.Item(2).Insert()
$wsHighPriority.Rows.Item(2).ClearFormats()

$wsHighPriority.Range("A2").Value = "Trend New"
$wsHighPriority.Range("B2").Value = 87
$wsHighPriority.Range("C2").Value = 84.5
$wsHighPriority.Range("D2").Value = 87
$wsHighPriority.Range("E2").Value = 86.09999999999999

$wsHighPriority.Range("A3").Value = "IUCN"
$wsHighPriority.Range("B3").Value = 16
$wsHighPriority.Range("C3").Value = 15.5
$wsHighPriority.Range("D3").Value = 14
$wsHighPriority.Range("E3").Value = 13.9

# Restore the originally active/selected sheet.
$wsTrends.Activate()
